$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 726, shifting existing rows 726:767 down to 727:768
$ws.Rows.Item(726).Insert()

# Populate the newly inserted row with the new data point (2026/01/26, Mon, hour 5, rank 165)
# Force column A to be stored as text (matching the rest of the date column) instead of
# letting Excel auto-convert the date-like string into a serial date value.
$ws.Cells.Item(726, 1).NumberFormat = "@"
$ws.Cells.Item(726, 1).Value = "2026/01/26"
$ws.Cells.Item(726, 1).ClearFormats()

$ws.Cells.Item(726, 2).Value = "月"
$ws.Cells.Item(726, 3).Value = 5
$ws.Cells.Item(726, 4).Value = 165
